# Apply the commit's changes to the workbook:
#  - Rename sheet "emi" to "source"
#  - Change its header cell A1 text from "SOURCES" to "SOURCE"
#  - Make "source" the active sheet/tab, with A2 selected
#    (this also clears the previous "grid" tab's active/selected state)

$wb = $excel.ActiveWorkbook

$wsSource = $wb.Worksheets.Item("emi")
$wsSource.Name = "source"
$wsSource.Range("A1").Value = "SOURCE"

$wsSource.Activate()
[void]$wsSource.Range("A2").Select()
